{"js": "// Updates the worksheet date header and the 25 division-problem\n// answer cells to the next day's generated values.\nconst replacements = [\n  [\"2025-08-10 Sunday\", \"2025-08-11 Monday\"],\n  [\"926\u00f76=154, 2\", \"964\u00f78=120, 4\"],\n  [\"557\u00f74=139, 1\", \"809\u00f78=101, 1\"],\n  [\"644\u00f77=92, 0\", \"858\u00f78=107, 2\"],\n  [\"265\u00f75=53, 0\", \"843\u00f76=140, 3\"],\n  [\"552\u00f79=61, 3\", \"651\u00f78=81, 3\"],\n  [\"398\u00f72=199, 0\", \"510\u00f78=63, 6\"],\n  [\"467\u00f79=51, 8\", \"481\u00f74=120, 1\"],\n  [\"738\u00f73=246, 0\", \"667\u00f79=74, 1\"],\n  [\"970\u00f79=107, 7\", \"635\u00f79=70, 5\"],\n  [\"661\u00f77=94, 3\", \"265\u00f76=44, 1\"],\n  [\"872\u00f76=145, 2\", \"523\u00f78=65, 3\"],\n  [\"691\u00f75=138, 1\", \"605\u00f78=75, 5\"],\n  [\"144\u00f75=28, 4\", \"855\u00f79=95, 0\"],\n  [\"100\u00f77=14, 2\", \"665\u00f74=166, 1\"],\n  [\"726\u00f79=80, 6\", \"662\u00f72=331, 0\"],\n  [\"911\u00f75=182, 1\", \"394\u00f72=197, 0\"],\n  [\"490\u00f73=163, 1\", \"848\u00f74=212, 0\"],\n  [\"985\u00f74=246, 1\", \"257\u00f75=51, 2\"],\n  [\"831\u00f79=92, 3\", \"663\u00f78=82, 7\"],\n  [\"652\u00f77=93, 1\", \"645\u00f74=161, 1\"],\n  [\"634\u00f79=70, 4\", \"201\u00f72=100, 1\"],\n  [\"799\u00f74=199, 3\", \"480\u00f78=60, 0\"],\n  [\"227\u00f72=113, 1\", \"433\u00f78=54, 1\"],\n  [\"733\u00f77=104, 5\", \"629\u00f79=69, 8\"],\n  [\"134\u00f77=19, 1\", \"226\u00f76=37, 4\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Updates the worksheet date header and the 25 division-problem\n# answer cells to the next day's generated values.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-08-10 Sunday\", \"2025-08-11 Monday\"),\n    @(\"926\u00f76=154, 2\", \"964\u00f78=120, 4\"),\n    @(\"557\u00f74=139, 1\", \"809\u00f78=101, 1\"),\n    @(\"644\u00f77=92, 0\", \"858\u00f78=107, 2\"),\n    @(\"265\u00f75=53, 0\", \"843\u00f76=140, 3\"),\n    @(\"552\u00f79=61, 3\", \"651\u00f78=81, 3\"),\n    @(\"398\u00f72=199, 0\", \"510\u00f78=63, 6\"),\n    @(\"467\u00f79=51, 8\", \"481\u00f74=120, 1\"),\n    @(\"738\u00f73=246, 0\", \"667\u00f79=74, 1\"),\n    @(\"970\u00f79=107, 7\", \"635\u00f79=70, 5\"),\n    @(\"661\u00f77=94, 3\", \"265\u00f76=44, 1\"),\n    @(\"872\u00f76=145, 2\", \"523\u00f78=65, 3\"),\n    @(\"691\u00f75=138, 1\", \"605\u00f78=75, 5\"),\n    @(\"144\u00f75=28, 4\", \"855\u00f79=95, 0\"),\n    @(\"100\u00f77=14, 2\", \"665\u00f74=166, 1\"),\n    @(\"726\u00f79=80, 6\", \"662\u00f72=331, 0\"),\n    @(\"911\u00f75=182, 1\", \"394\u00f72=197, 0\"),\n    @(\"490\u00f73=163, 1\", \"848\u00f74=212, 0\"),\n    @(\"985\u00f74=246, 1\", \"257\u00f75=51, 2\"),\n    @(\"831\u00f79=92, 3\", \"663\u00f78=82, 7\"),\n    @(\"652\u00f77=93, 1\", \"645\u00f74=161, 1\"),\n    @(\"634\u00f79=70, 4\", \"201\u00f72=100, 1\"),\n    @(\"799\u00f74=199, 3\", \"480\u00f78=60, 0\"),\n    @(\"227\u00f72=113, 1\", \"433\u00f78=54, 1\"),\n    @(\"733\u00f77=104, 5\", \"629\u00f79=69, 8\"),\n    @(\"134\u00f77=19, 1\", \"226\u00f76=37, 4\"),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # WdFindWrap.wdFindContinue = 1, WdReplace.wdReplaceAll = 2\n    $find.Execute($pair[0], $true, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
